# Insert a new weekly price record as row 46 in the "Espinaca" price
# history sheet, pushing the existing rows 46-82 down to 47-83.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 46 (shifts rows 46..82 down to 47..83).
$ws.Rows.Item(46).Insert()

# Populate the new row 46 with the latest weekly record.
$ws.Cells.Item(46, 1).Value = 4
$ws.Cells.Item(46, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(46, 3).Value = "Los Lagos"
$ws.Cells.Item(46, 4).Value = 45167
$ws.Cells.Item(46, 5).Value = 10
$ws.Cells.Item(46, 6).Value = 100112012
$ws.Cells.Item(46, 7).Value = "Espinaca"
$ws.Cells.Item(46, 8).Value = "Sin especificar"
$ws.Cells.Item(46, 9).Value = "Primera"
$ws.Cells.Item(46, 10).Value = 35
$ws.Cells.Item(46, 11).Value = 13000
$ws.Cells.Item(46, 12).Value = 13000
$ws.Cells.Item(46, 13).Value = 13000
$ws.Cells.Item(46, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(46, 15).Value = "Región Metropolitana"
$ws.Cells.Item(46, 16).Value = 1300
$ws.Cells.Item(46, 17).Value = 10
$ws.Cells.Item(46, 18).Value = "Hortaliza"
